$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 74820.664
$ws.Range("I6").Value = 500.33334
$ws.Range("J6").Value = 149141
$ws.Range("K6").Value = 1501.00002
$ws.Range("L6").Value = 447423
$ws.Range("M6").Value = -1389.00002
$ws.Range("N6").Value = -447647
$ws.Range("H15").Value = 3063.1667
$ws.Range("I15").Value = 3063.1667
$ws.Range("K15").Value = 9189.500100000001
$ws.Range("M15").Value = -9020.500100000001
$ws.Range("H58").Value = 209.61539
$ws.Range("I58").Value = 190.625
$ws.Range("J58").Value = 240
$ws.Range("K58").Value = 571.875
$ws.Range("L58").Value = 720
$ws.Range("M58").Value = -421.875
$ws.Range("N58").Value = -1020
$ws.Range("H138").Value = 2453043
$ws.Range("I138").Value = 1571.05
$ws.Range("J138").Value = 5955145.5
$ws.Range("K138").Value = 4713.15
$ws.Range("L138").Value = 17865436.5
$ws.Range("M138").Value = 426.8500000000004
$ws.Range("N138").Value = -17875716.5
$ws.Range("H141").Value = 2690.5
$ws.Range("I141").Value = 1420.2858
$ws.Range("K141").Value = 4260.857400000001
$ws.Range("M141").Value = 919.1425999999992
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3413.19
$ws.Range("I32").Value = 2842.8354
$ws.Range("J32").Value = 5558.8096
$ws.Range("K32").Value = 2842.8354
$ws.Range("L32").Value = 5558.8096
$ws.Range("M32").Value = -2555.8354
$ws.Range("N32").Value = -6132.8096
$ws.Range("H121").Value = 100000
$ws.Range("J121").Value = 100000
$ws.Range("L121").Value = 100000
$ws.Range("N121").Value = -103494
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13995.167
$ws.Range("I26").Value = 13995.167
$ws.Range("K26").Value = 13995.167
$ws.Range("M26").Value = -13703.167
$ws.Range("H99").Value = 1046.3636
$ws.Range("I99").Value = 906.125
$ws.Range("J99").Value = 1420.3334
$ws.Range("K99").Value = 906.125
$ws.Range("L99").Value = 1420.3334
$ws.Range("M99").Value = 591.875
$ws.Range("N99").Value = -4416.3334
$ws.Range("H105").Value = 17243626
$ws.Range("I105").Value = 27779912
$ws.Range("J105").Value = 2427.2727
$ws.Range("K105").Value = 27779912
$ws.Range("L105").Value = 2427.2727
$ws.Range("M105").Value = -27778165
$ws.Range("N105").Value = -5921.2727
$ws.Range("H122").Value = 53314.285
$ws.Range("J122").Value = 53314.285
$ws.Range("L122").Value = 53314.285
$ws.Range("N122").Value = -63114.285
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 50780
$ws.Range("J130").Value = 50780
$ws.Range("L130").Value = 50780
$ws.Range("N130").Value = -60820
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 51296.668
$ws.Range("J98").Value = 51296.668
$ws.Range("L98").Value = 51296.668
$ws.Range("N98").Value = -55788.668
$ws.Range("H105").Value = 610
$ws.Range("I105").Value = 610
$ws.Range("K105").Value = 610
$ws.Range("M105").Value = 1137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 169908.66
$ws.Range("I56").Value = 169908.66
$ws.Range("K56").Value = 169908.66
$ws.Range("M56").Value = -169378.66
$ws.Range("H75").Value = 1392.1666
$ws.Range("I75").Value = 1392.1666
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 4176.4998
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -3178.4998
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1392.1666
$ws.Range("I78").Value = 1392.1666
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 12529.4994
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -7537.499400000001
$ws.Range("N78").ClearContents()
$ws.Range("H116").Value = 122704.63
$ws.Range("J116").Value = 85207.75
$ws.Range("L116").Value = 255623.25
$ws.Range("N116").Value = -262507.25
$ws.Range("H122").Value = 594.175
$ws.Range("J122").Value = 876.1905
$ws.Range("L122").Value = 7885.7145
$ws.Range("N122").Value = -12785.7145
$ws.Range("H131").Value = 1024.0883
$ws.Range("I131").Value = 387
$ws.Range("J131").Value = 1133.931
$ws.Range("K131").Value = 1161
$ws.Range("L131").Value = 3401.793
$ws.Range("M131").Value = 3879
$ws.Range("N131").Value = -13481.793
$ws.Range("H132").Value = 2414.1738
$ws.Range("I132").Value = 1795.125
$ws.Range("J132").Value = 2744.3333
$ws.Range("K132").Value = 16156.125
$ws.Range("L132").Value = 24698.9997
$ws.Range("M132").Value = -13626.125
$ws.Range("N132").Value = -29758.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 35728.555
$ws.Range("J34").Value = 35728.555
$ws.Range("L34").Value = 35728.555
$ws.Range("N34").Value = -36264.555
$ws.Range("H57").Value = 18490
$ws.Range("J57").Value = 18390
$ws.Range("L57").Value = 18390
$ws.Range("N57").Value = -20030
$ws.Range("H76").Value = 35728.555
$ws.Range("J76").Value = 35728.555
$ws.Range("L76").Value = 35728.555
$ws.Range("N76").Value = -36358.555
$ws.Range("H79").Value = 35728.555
$ws.Range("J79").Value = 35728.555
$ws.Range("L79").Value = 35728.555
$ws.Range("N79").Value = -37912.555
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2043.1111
$ws.Range("I82").Value = 1750.25
$ws.Range("J82").Value = 2277.4
$ws.Range("K82").Value = 1750.25
$ws.Range("L82").Value = 2277.4
$ws.Range("M82").Value = -1389.25
$ws.Range("N82").Value = -2999.4
$ws.Range("H85").Value = 2043.1111
$ws.Range("I85").Value = 1750.25
$ws.Range("J85").Value = 2277.4
$ws.Range("K85").Value = 1750.25
$ws.Range("L85").Value = 2277.4
$ws.Range("M85").Value = -502.25
$ws.Range("N85").Value = -4773.4
$ws.Range("H100").Value = 1607
$ws.Range("I100").Value = 1415.6154
$ws.Range("J100").Value = 1762.5
$ws.Range("K100").Value = 1415.6154
$ws.Range("L100").Value = 1762.5
$ws.Range("M100").Value = -874.6153999999999
$ws.Range("N100").Value = -2844.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 49448
$ws.Range("J95").Value = 49448
$ws.Range("L95").Value = 49448
$ws.Range("N95").Value = -54940
$ws.Range("H136").Value = 33095.47
$ws.Range("I136").Value = 30496.03
$ws.Range("J136").Value = 36041.5
$ws.Range("K136").Value = 91488.09
$ws.Range("L136").Value = 108124.5
$ws.Range("M136").Value = -88938.09
$ws.Range("N136").Value = -113224.5
